$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-strings table was reordered by a re-run of text processing.
# Effectively, the text shown in these specific rows (column A) changed
# while the row's count (column B) stayed put. Update each affected cell
# to its new text value.
$ws.Range("A19").Value = "небогатый товар"
$ws.Range("A20").Value = "крамными товар"

$ws.Range("A28").Value = "суровский товар"
$ws.Range("A30").Value = "медный товар"

$ws.Range("A31").Value = "внутренний товар"
$ws.Range("A32").Value = "питейный припасы"

$ws.Range("A35").Value = "произрастание"
$ws.Range("A37").Value = "купецкий товар"
$ws.Range("A38").Value = "заморский товар"

$ws.Range("A39").Value = "меховой товар"
$ws.Range("A41").Value = "надлежащий товар"
$ws.Range("A42").Value = "рукодельный товар"
$ws.Range("A43").Value = "домовый товар"
